$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Reset the selection on Sheet1 to the full data range (no more "C4" leftover
# selection from the old file) before we add/activate the new sheet.
$sheet1.Range("A1:C3").Select() | Out-Null

# Add the new sheet right after Sheet1 and name it.
$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "New Sheet"

# Populate the new sheet's contents.
$newSheet.Range("A1").Value = "Other stuff"

$newSheet.Range("I6").Value = "name"
$newSheet.Range("J6").Value = "x_origin"
$newSheet.Range("K6").Value = "y_origin"

$newSheet.Range("I7").Value = "Zone 3"
$newSheet.Range("J7").Value = 0
$newSheet.Range("K7").Value = 0

$newSheet.Range("I8").Value = "Zone 4"
$newSheet.Range("J8").Value = 5
$newSheet.Range("K8").Value = 10

$newSheet.Columns.Item(1).AutoFit() | Out-Null

# Match the authored selection/active-cell state on the new sheet.
$newSheet.Range("I6").Select() | Out-Null
